$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("registro")

# Update header row to lowercase labels
$ws.Range("A1").Value = "correo electronico"
$ws.Range("B1").Value = "nombre"
$ws.Range("C1").Value = "contrasena"
$ws.Range("D1").Value = "nombre_tablero"

# Update the active selection
$ws.Activate()
$ws.Range("F6").Select()
